$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add the new row 31 with the mail log entry
$ws.Range("A31").Value = "Demo inplannen"
$ws.Range("B31").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C31").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D31").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E31").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F31").Value = "2025-08-14 21:49:21"
$ws.Range("G31").Value = "Nee"
$ws.Range("H31").Value = "Ja"
$ws.Range("I31").Value = "Nee"
$ws.Range("J31").Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$ranges = @("D2:D30", "G2:G30", "H2:H30", "I2:I30", "J2:J30")
$newRanges = @("D2:D31", "G2:G31", "H2:H31", "I2:I31", "J2:J31")

for ($k = 0; $k -lt $ranges.Count; $k++) {
    $fcs = $ws.Range($ranges[$k]).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($newRanges[$k]))
    }
}

# Update the Dashboard count for "Intern verzoek / Actie voor medewerker"
$dash.Range("B2").Value = 23
